$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 590, shifting existing rows 590:659 down to 591:660
$ws.Rows.Item(590).Insert()

# Populate the newly inserted row 590 with the new data record
$ws.Range("A590").Value = 10
$ws.Range("B590").Value = "Vega Modelo de Temuco"
$ws.Range("C590").Value = "La Araucanía"
$ws.Range("D590").Value = 44918
$ws.Range("E590").Value = 9
$ws.Range("F590").Value = 100112032
$ws.Range("G590").Value = "Zapallo italiano"
$ws.Range("H590").Value = "Sin especificar"
$ws.Range("I590").Value = "Primera"
$ws.Range("J590").Value = 210
$ws.Range("K590").Value = 9000
$ws.Range("L590").Value = 10000
$ws.Range("M590").Value = 9405
$ws.Range("N590").Value = "$/caja 50 unidades"
$ws.Range("O590").Value = "Región del Maule"
$ws.Range("P590").Value = 188
$ws.Range("Q590").Value = 50
$ws.Range("R590").Value = "Hortaliza"
